$d = $word.ActiveDocument

# --- Change 1: Table of Contents entry "4. Feature Specifications" gets a
#     line break followed by a new bullet line describing the Factorial
#     operation. Scope the Find/Replace to the single TOC paragraph (the
#     one using the "ListNumber" style) so the later "4. Feature
#     Specifications" Heading1 occurrence (with lastRenderedPageBreak) is
#     left untouched.
$tocPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "4. Feature Specifications" `
        -and $para.Style.NameLocal -eq "List Number") {
        $tocPara = $para
        break
    }
}

if ($tocPara -eq $null) {
    throw "Could not locate the 'List Number' TOC paragraph for '4. Feature Specifications'"
}

$r = $tocPara.Range
$r.Find.Execute("4. Feature Specifications", $true, $false, $false, $false, $false, $true, 1, $false, `
    "4. Feature Specifications^l• Factorial Operation: Calculates the factorial of a number", 2)

# --- Change 2: add a new "Factorial" / "factorial" row to the end of the
#     Option/Function table (the first table in the document).
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Factorial"
$newRow.Cells.Item(2).Range.Text = "factorial"
